$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.333.84"
$ws.Range("E2").Value = "  +0.45%  "
$ws.Range("D3").Value = "1.668.72"
$ws.Range("E3").Value = "  +0.79%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'220.46"
$ws.Range("E5").Value = "  +1.02%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D8").Value = "'0.2651"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").Value = "'0.06372"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").Value = "'20.90"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").Value = "'0.07828"
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("D12").Value = "'4.529"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "1.673.93"
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("D14").Value = "1.897.22"
$ws.Range("E14").Value = "  +0.76%  "
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").Value = "0.0₅8171"
$ws.Range("D17").Value = "'65.95"
$ws.Range("E17").Value = "  +0.85%  "
$ws.Range("D18").Value = "26.336.87"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("D21").Value = "'198.57"
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").Value = "'6.053"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'146.38"
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").Value = "'16.21"
$ws.Range("E28").Value = "  +1.32%  "
$ws.Range("D29").Value = "'1.504"
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("D30").Value = "'0.05907"
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").Value = "'1.285"
$ws.Range("E31").Value = "  +0.61%  "
$ws.Range("D32").Value = "'3.553"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("D33").Value = "'3.325"
$ws.Range("E33").Value = "  +1.62%  "
$ws.Range("D34").Value = "'1.605"
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Value = "'2.836"
$ws.Range("E35").Value = "  +0.73%  "
$ws.Range("D36").Value = "'0.9623"
$ws.Range("E36").Value = "  +1.04%  "
$ws.Range("D37").Value = "'2.433"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("D38").Value = "'0.5825"
$ws.Range("E38").Value = "  +0.98%  "
$ws.Range("D39").Value = "'0.01618"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "'5.958"
$ws.Range("D41").Value = "1.077.82"
$ws.Range("E41").Value = "  +2.93%  "
$ws.Range("D42").Value = "'0.8596"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").Value = "'102.94"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").Value = "1.807.71"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "'58.64"
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("D48").Value = "'1.013"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").Value = "'8.051"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'0.05152"
